$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44432
$ws.Range("M2").Value = 20

$ws.Range("D3").Value = 44435
$ws.Range("M3").Value = 40
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 2000

$ws.Range("D4").Value = 44517
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 27000
$ws.Range("O4").Value = 27000
$ws.Range("P4").Value = 27000
$ws.Range("S4").Value = 2700

$ws.Range("D5").Value = 44517
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 25000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 25000
$ws.Range("S5").Value = 2500

$ws.Range("D6").Value = 44476
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("S6").Value = 2000

$ws.Range("D7").Value = 44473
$ws.Range("M7").Value = 180
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("S7").Value = 2000

$ws.Range("D8").Value = 44466
$ws.Range("M8").Value = 60

$ws.Range("D10").Value = 44503
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 30000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 30000
$ws.Range("S10").Value = 3000

$ws.Range("D11").Value = 44503
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 25000
$ws.Range("S11").Value = 2500
